# Insert a new weekly price observation as row 191 (pushing the existing
# rows 191-240 down to 192-241), matching the logic of a daily/weekly
# "Fruta / hortaliza" price update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 191..240 down to 192..241, creating an empty row 191.
$ws.Rows.Item(191).Insert()

# Populate the newly inserted row 191 with the new observation.
$ws.Cells.Item(191, 1).Value = 11
$ws.Cells.Item(191, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(191, 3).Value = "Bíobío"
$ws.Cells.Item(191, 4).Value = 44932
$ws.Cells.Item(191, 5).Value = 8
$ws.Cells.Item(191, 6).Value = 100112003
$ws.Cells.Item(191, 7).Value = "Ajo"
$ws.Cells.Item(191, 8).Value = "Chino"
$ws.Cells.Item(191, 9).Value = "Primera"
$ws.Cells.Item(191, 10).Value = 310
$ws.Cells.Item(191, 11).Value = 13000
$ws.Cells.Item(191, 12).Value = 14000
$ws.Cells.Item(191, 13).Value = 13516
$ws.Cells.Item(191, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(191, 15).Value = "China"
$ws.Cells.Item(191, 16).Value = 1352
$ws.Cells.Item(191, 17).Value = 10
$ws.Cells.Item(191, 18).Value = "Hortaliza"
